# cpu_upl/controls.xlsx — add a new control-signal column "ctrl_reg_input_mux"
# (column P) to the truth table on Sheet1, and update the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header in P1, bold like the other header cells (style index 2 in the
# original file: fontId=1 / bold, no custom number format).
$ws.Range("P1").Value = "ctrl_reg_input_mux"
$ws.Range("P1").Font.Bold = $true

# New column values, row by row (row 16 is the only "1"; everything else is 0).
$pValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 0
}

foreach ($row in $pValues.Keys) {
    $ws.Cells.Item($row, 16).Value = $pValues[$row]
}

# Reflect the author's final selection/scroll position on the sheet
# (scrolled one column right, so column B is leftmost/visible, with O22
# as the active cell).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("O22").Select()
